$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2 ("VR Template-Interactobot_rand"): add Run 1..Run 5 columns + data ---

# Header row: C1..G1 = "Run 1".."Run 5"
$ws2.Range("C1").Value = "Run 1"
$ws2.Range("D1").Value = "Run 2"
$ws2.Range("E1").Value = "Run 3"
$ws2.Range("F1").Value = "Run 4"
$ws2.Range("G1").Value = "Run 5"

# Time column A2:A11 = 30,60,...,300 (sheet2 grows from 5 rows to 10, matching sheet1's time series)
$times = @(30, 60, 90, 120, 150, 180, 210, 240, 270, 300)
for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $times[$i]
}

# Per-run counts for columns C..G (Run1..Run5), rows 2..11
$runData = @(
    @(0,2,0,1,0),
    @(0,2,0,1,0),
    @(0,2,1,2,0),
    @(0,3,1,2,0),
    @(0,3,1,2,0),
    @(0,3,2,2,1),
    @(0,3,2,2,1),
    @(0,3,2,3,3),
    @(0,3,2,3,3),
    @(0,3,2,3,3)
)

for ($i = 0; $i -lt $runData.Length; $i++) {
    $row = $i + 2
    $vals = $runData[$i]
    $ws2.Cells.Item($row, 3).Value = $vals[0]
    $ws2.Cells.Item($row, 4).Value = $vals[1]
    $ws2.Cells.Item($row, 5).Value = $vals[2]
    $ws2.Cells.Item($row, 6).Value = $vals[3]
    $ws2.Cells.Item($row, 7).Value = $vals[4]
}

# Coverage formula in column B, rows 2..11: average coverage across the 5 runs
for ($row = 2; $row -le 11; $row++) {
    $ws2.Cells.Item($row, 2).Formula = "=SUM(C" + $row + ":G" + $row + ")/5/11"
}

# Re-apply the 2-decimal numeric format to the whole Coverage column (B) so every
# cell (old + newly-added rows) shares one consistent style, and so the column
# picks up the same default-width column record sheet1 already carries.
$ws2.Columns.Item(2).ColumnWidth = 10.83203125
$ws2.Columns.Item(2).NumberFormat = "0.00"

# Page setup (portrait, paper size 9 = A4) now recorded for sheet2, matching sheet1's
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Sheet2 becomes the active/visible tab with its own selection
$ws2.Range("E16").Select() | Out-Null

# --- Sheet1 ("VR Template-Interactobot"): selection narrows from A2:B11 to A2:A11 ---
$ws1.Range("A2:A11").Select() | Out-Null

# Sheet2 is the tab shown/active when the workbook is reopened
$ws2.Activate() | Out-Null

Write-Output "done"
